# Revert "Adding Unit 1 homework draft":
#  - Remove the Table1 ListObject on Sheet1 (convert back to a plain range)
#  - Shift the Date column (Sheet1 column E) back by 980 days to the
#    original 2016 dates (matching the still-unchanged pivot cache)
#  - Remove the two extra homework-draft note rows added at the bottom of
#    Sheet2 (the shared strings that went with them fall out naturally)
#  - Restore the original active-cell selection per sheet

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the table, keep the data/autofilter range -------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$table1 = $ws1.ListObjects.Item(1)
$null = $table1.Unlist()

# --- Sheet1: shift every Date value (column E, rows 2-214) back 980 days
for ($r = 2; $r -le 214; $r++) {
    $cell = $ws1.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 - 980
}

# --- Sheet2: remove the two note rows (14 and 15) that were appended ----
$ws2 = $wb.Worksheets.Item("Sheet2")
$null = $ws2.Range("A14:A15").EntireRow.Delete()

# --- Restore the per-sheet active-cell selections ------------------------
$ws1.Activate()
$null = $ws1.Range("A2").Select()

$ws2.Activate()
$null = $ws2.Range("D4").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$null = $ws3.Range("A10").Select()

$ws1.Activate()
